$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-running the notebook on updated data reshuffled the tie-breaking order
# of words that share the same count in the frequency table. Counts (column B)
# are unchanged; only the word (column A) assigned to a handful of rows moves.
$ws.Range("A16").Value = "колеса"
$ws.Range("A17").Value = "полотно"
$ws.Range("A18").Value = "Крымскую соль"
$ws.Range("A19").Value = "говядина"
$ws.Range("A21").Value = "сено"
$ws.Range("A24").Value = "сахар"
$ws.Range("A25").Value = "выбойка"
$ws.Range("A26").Value = "чулок"
$ws.Range("A27").Value = "шелк"
$ws.Range("A35").Value = "горшок"
$ws.Range("A36").Value = "конь"
$ws.Range("A37").Value = "обод"
$ws.Range("A38").Value = "веревка"
$ws.Range("A39").Value = "ром"
$ws.Range("A40").Value = "гвоздь"
$ws.Range("A41").Value = "рогожа"
$ws.Range("A42").Value = "овца"
$ws.Range("A43").Value = "замок"
$ws.Range("A45").Value = "дуга"
$ws.Range("A46").Value = "покроми"
$ws.Range("A47").Value = "роза"
$ws.Range("A48").Value = "хомут"
$ws.Range("A49").Value = "гумми"
$ws.Range("A51").Value = "сосуд"
$ws.Range("A52").Value = "нитка"
$ws.Range("A53").Value = "брусья"
$ws.Range("A54").Value = "котел"
$ws.Range("A55").Value = "бечева"
$ws.Range("A56").Value = "скотский кожа"
